# Rename sheet 'Data' to 'Data table' and make it the active sheet/tab
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Name = "Data table"
$ws.Activate()
